# ST - up to final boss fight
# Fills in the boss-4 / boss-5 fight splits (rows 197-210), the gap row 211
# (raw split time only), shifts the old "end level" row down to 212 and
# adds a new "start fight" row 213, leaves a gap (214-215) and shifts the
# old "boss fight end (white screen)" row down to 216.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V3")

function Set-SplitRow($Row, $Place, $B, $C) {
    if ($null -ne $Place) { $ws.Cells.Item($Row, 1).Value = $Place }
    if ($null -ne $B) { $ws.Cells.Item($Row, 2).Value = $B }
    if ($null -ne $C) { $ws.Cells.Item($Row, 3).Value = $C }
    $ws.Range("D$Row").Formula = "=IF(B$Row=`"`",`"-`",IF(C$Row=`"`",`"-`",B$Row-C$Row))"
}

Set-SplitRow 197 "Boss 3 HP = 0"  165246 153517
Set-SplitRow 198 "Boss 4 appears" 165462 153699
Set-SplitRow 199 "Boss 4 HP = 40" 165547 153784
Set-SplitRow 200 "Boss 4 HP = 30" 165788 154025
Set-SplitRow 201 "Boss 4 HP = 20" 166029 154266
Set-SplitRow 202 "Boss 4 HP = 10" 166282 154507
Set-SplitRow 203 "Boss 4 HP = 0"  166526 154747
Set-SplitRow 204 "Boss 5 appears" 167321 155542
Set-SplitRow 205 "Boss 5 HP = 50" 167361 155582
Set-SplitRow 206 "Boss 5 HP = 40" 167603 155822
Set-SplitRow 207 "Boss 5 HP = 30" 167849 156062
Set-SplitRow 208 "Boss 5 HP = 20" 168155 156302
Set-SplitRow 209 "Boss 5 HP = 10" 168559 156542
Set-SplitRow 210 "Boss 5 HP = 00" 168967 156782

# Row 211 only ever got a raw C split value logged - no place label, B or D.
$ws.Cells.Item(211, 3).Value = 157395

# The old "end level" row (previously row 199) moves down to 212, now with
# a C value added.
Set-SplitRow 212 "end level" 175028 162843

# New "start fight" row for the final boss.
Set-SplitRow 213 "start fight" 175962 163774

# Rows 214-215 stay empty (gap before the last logged row).

# The old "boss fight end (white screen)" row (previously row 200) moves
# down to 216; it never got a C value, so D evaluates to "-".
$ws.Cells.Item(216, 1).Value = "boss fight end (white screen)"
$ws.Cells.Item(216, 2).Value = 179257
$ws.Range("D216").Formula = "=IF(B216=`"`",`"-`",IF(C216=`"`",`"-`",B216-C216))"

# Move the view down to where the new rows were entered, matching the
# author's on-screen state after the edit.
$ws.Cells.Item(214, 2).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 202
